$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: split "deviationPercent" into separate yielding/hardening columns,
# and shift "optimizerName" header from I1 into J1.
$ws.Range("H1").Value = "yielding_deviationPercent"
$ws.Range("I1").Value = "hardening_deviationPercent"
$ws.Range("J1").Value = "optimizerName"
# I1 now mirrors H1's numeric-style header formatting (#,##0, centered, wrap).
$ws.Range("I1").NumberFormat = $ws.Range("H1").NumberFormat

# Data row 2: update material/hardeningLaw values and populate the new
# yielding/hardening deviation + optimizerName columns.
$ws.Range("D2").Value = "DP1000_room"
$ws.Range("E2").Value = "Swift"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 2
# I2 becomes a numeric cell like H2, so pick up the same #,##0 formatting.
$ws.Range("I2").NumberFormat = $ws.Range("H2").NumberFormat
$ws.Range("J2").Value = "BO"
